# Weekly refresh of "Hortaliza, Vega Central Mapocho de Santiago - Oregano" price series.
# Each weekly run re-pulls the source feed: existing rows keep their identity
# (Mercado/Categoria/Calidad/unit columns, A/B/C/E/F/G/H/I/N/O/Q/R) but the
# date (D), volume (J) and price columns (K/L/M/P) are refreshed in place,
# and the newest week's observation is appended as a new row at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44503
$ws.Cells.Item(2, 10).Value = 16
$ws.Cells.Item(2, 11).Value = 8000
$ws.Cells.Item(2, 12).Value = 9000
$ws.Cells.Item(2, 13).Value = 8500
$ws.Cells.Item(2, 16).Value = 2833
$ws.Cells.Item(3, 4).Value = 44461
$ws.Cells.Item(3, 10).Value = 16
$ws.Cells.Item(3, 11).Value = 9500
$ws.Cells.Item(3, 13).Value = 9750
$ws.Cells.Item(3, 16).Value = 3250
$ws.Cells.Item(4, 4).Value = 44370
$ws.Cells.Item(4, 11).Value = 10000
$ws.Cells.Item(4, 12).Value = 10500
$ws.Cells.Item(4, 13).Value = 10250
$ws.Cells.Item(4, 16).Value = 3417
$ws.Cells.Item(5, 4).Value = 44377
$ws.Cells.Item(5, 10).Value = 16
$ws.Cells.Item(5, 12).Value = 10500
$ws.Cells.Item(5, 13).Value = 10250
$ws.Cells.Item(5, 16).Value = 3417
$ws.Cells.Item(6, 4).Value = 44314
$ws.Cells.Item(7, 4).Value = 44412
$ws.Cells.Item(7, 10).Value = 25
$ws.Cells.Item(7, 11).Value = 10000
$ws.Cells.Item(7, 12).Value = 10500
$ws.Cells.Item(7, 13).Value = 10260
$ws.Cells.Item(7, 16).Value = 3420
$ws.Cells.Item(8, 4).Value = 44266
$ws.Cells.Item(8, 10).Value = 160
$ws.Cells.Item(9, 4).Value = 44433
$ws.Cells.Item(9, 12).Value = 10500
$ws.Cells.Item(9, 13).Value = 10250
$ws.Cells.Item(9, 16).Value = 3417
$ws.Cells.Item(10, 4).Value = 44335
$ws.Cells.Item(10, 12).Value = 10000
$ws.Cells.Item(10, 13).Value = 10000
$ws.Cells.Item(10, 16).Value = 3333
$ws.Cells.Item(12, 4).Value = 44482
$ws.Cells.Item(12, 11).Value = 9000
$ws.Cells.Item(12, 13).Value = 9500
$ws.Cells.Item(12, 16).Value = 3167
$ws.Cells.Item(13, 4).Value = 44405
$ws.Cells.Item(13, 12).Value = 10500
$ws.Cells.Item(13, 13).Value = 10250
$ws.Cells.Item(13, 16).Value = 3417
$ws.Cells.Item(14, 4).Value = 44435
$ws.Cells.Item(15, 4).Value = 44328
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 10000
$ws.Cells.Item(15, 16).Value = 3333
$ws.Cells.Item(16, 4).Value = 44175
$ws.Cells.Item(16, 10).Value = 70
$ws.Cells.Item(16, 11).Value = 12000
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).Value = 12000
$ws.Cells.Item(16, 16).Value = 4000
$ws.Cells.Item(17, 4).Value = 44475
$ws.Cells.Item(17, 10).Value = 16
$ws.Cells.Item(17, 11).Value = 9000
$ws.Cells.Item(17, 13).Value = 9500
$ws.Cells.Item(17, 16).Value = 3167
$ws.Cells.Item(18, 4).Value = 44419
$ws.Cells.Item(18, 10).Value = 16
$ws.Cells.Item(18, 12).Value = 10000
$ws.Cells.Item(18, 13).Value = 10000
$ws.Cells.Item(18, 16).Value = 3333
$ws.Cells.Item(19, 4).Value = 44293
$ws.Cells.Item(19, 12).Value = 10000
$ws.Cells.Item(19, 13).Value = 10000
$ws.Cells.Item(19, 16).Value = 3333
$ws.Cells.Item(20, 4).Value = 44398
$ws.Cells.Item(20, 12).Value = 10500
$ws.Cells.Item(20, 13).Value = 10250
$ws.Cells.Item(20, 16).Value = 3417
$ws.Cells.Item(21, 4).Value = 44321
$ws.Cells.Item(21, 10).Value = 25
$ws.Cells.Item(21, 12).Value = 10000
$ws.Cells.Item(21, 13).Value = 10000
$ws.Cells.Item(21, 16).Value = 3333
$ws.Cells.Item(22, 4).Value = 44349
$ws.Cells.Item(22, 10).Value = 12
$ws.Cells.Item(23, 4).Value = 44454
$ws.Cells.Item(23, 10).Value = 16
$ws.Cells.Item(23, 11).Value = 9500
$ws.Cells.Item(23, 13).Value = 9750
$ws.Cells.Item(23, 16).Value = 3250
$ws.Cells.Item(24, 4).Value = 44300
$ws.Cells.Item(24, 12).Value = 10000
$ws.Cells.Item(24, 13).Value = 10000
$ws.Cells.Item(24, 16).Value = 3333
$ws.Cells.Item(25, 4).Value = 44356
$ws.Cells.Item(25, 11).Value = 10000
$ws.Cells.Item(25, 13).Value = 10000
$ws.Cells.Item(25, 16).Value = 3333
$ws.Cells.Item(26, 4).Value = 44342
$ws.Cells.Item(26, 10).Value = 17
$ws.Cells.Item(26, 11).Value = 10000
$ws.Cells.Item(26, 13).Value = 10000
$ws.Cells.Item(26, 16).Value = 3333
$ws.Cells.Item(27, 4).Value = 44363
$ws.Cells.Item(28, 4).Value = 44195
$ws.Cells.Item(28, 10).Value = 30
$ws.Cells.Item(28, 11).Value = 10000
$ws.Cells.Item(28, 13).Value = 10000
$ws.Cells.Item(28, 16).Value = 3333
$ws.Cells.Item(29, 4).Value = 44426
$ws.Cells.Item(29, 12).Value = 10500
$ws.Cells.Item(29, 13).Value = 10250
$ws.Cells.Item(29, 16).Value = 3417
$ws.Cells.Item(30, 4).Value = 44279
$ws.Cells.Item(32, 4).Value = 44272
$ws.Cells.Item(32, 10).Value = 70
$ws.Cells.Item(32, 12).Value = 10000
$ws.Cells.Item(32, 13).Value = 10000
$ws.Cells.Item(32, 16).Value = 3333
$ws.Cells.Item(33, 4).Value = 44447
$ws.Cells.Item(33, 10).Value = 16
$ws.Cells.Item(33, 11).Value = 10000
$ws.Cells.Item(33, 12).Value = 10500
$ws.Cells.Item(33, 13).Value = 10250
$ws.Cells.Item(33, 16).Value = 3417
$ws.Cells.Item(34, 4).Value = 44510
$ws.Cells.Item(34, 10).Value = 16
$ws.Cells.Item(34, 11).Value = 9000
$ws.Cells.Item(34, 12).Value = 10000
$ws.Cells.Item(34, 13).Value = 9500
$ws.Cells.Item(34, 16).Value = 3167
$ws.Cells.Item(35, 4).Value = 44468
$ws.Cells.Item(35, 10).Value = 16
$ws.Cells.Item(35, 12).Value = 11000
$ws.Cells.Item(35, 13).Value = 10500
$ws.Cells.Item(35, 16).Value = 3500
$ws.Cells.Item(36, 4).Value = 44517
$ws.Cells.Item(36, 11).Value = 9000
$ws.Cells.Item(36, 13).Value = 9500
$ws.Cells.Item(36, 16).Value = 3167
$ws.Cells.Item(37, 4).Value = 44391
$ws.Cells.Item(37, 10).Value = 16
$ws.Cells.Item(38, 4).Value = 44181
$ws.Cells.Item(38, 10).Value = 10
$ws.Cells.Item(38, 12).Value = 12000
$ws.Cells.Item(38, 13).Value = 11000
$ws.Cells.Item(38, 16).Value = 3667

# Append the newest weekly observation as a new row at the bottom of the table.
$newRow = 39
$ws.Cells.Item($newRow, 1).Value = 9
$ws.Cells.Item($newRow, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($newRow, 3).Value = "Metropolitana"
$ws.Cells.Item($newRow, 4).Value = 44307
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value = 13
$ws.Cells.Item($newRow, 6).Value = 100112029
$ws.Cells.Item($newRow, 7).Value = "Orégano"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 160
$ws.Cells.Item($newRow, 11).Value = 10000
$ws.Cells.Item($newRow, 12).Value = 10000
$ws.Cells.Item($newRow, 13).Value = 10000
$ws.Cells.Item($newRow, 14).Value = "$/docena de atados"
$ws.Cells.Item($newRow, 15).Value = "Región Metropolitana"
$ws.Cells.Item($newRow, 16).Value = 3333
$ws.Cells.Item($newRow, 17).Value = 3
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
